$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("K3").Value2 = 19
$ws.Range("L3").Value2 = 1.13
$ws.Range("M3").Value2 = 6
$ws.Range("N3").Value2 = 1.44
$ws.Range("O3").Value2 = 2.7
# Row 4
$ws.Range("G4").Value2 = 5.25
$ws.Range("I4").Value2 = 1.62
$ws.Range("J4").Value2 = 1.04
$ws.Range("K4").Value2 = 12
$ws.Range("AB4").Value2 = 17
# Row 5
$ws.Range("G5").Value2 = 2.63
$ws.Range("I5").Value2 = 2.7
$ws.Range("U5").Value2 = 13
$ws.Range("W5").Value2 = 26
$ws.Range("Z5").Value2 = 10
$ws.Range("AE5").Value2 = 8.5
$ws.Range("AG5").Value2 = 10
$ws.Range("AH5").Value2 = 26
$ws.Range("AI5").Value2 = 21
$ws.Range("AJ5").Value2 = 29
# Row 7
$ws.Range("N7").Value2 = 2.6
$ws.Range("O7").Value2 = 1.48
# Row 14
$ws.Range("G14").Value2 = 5
$ws.Range("H14").Value2 = 3.75
$ws.Range("I14").Value2 = 1.62
$ws.Range("J14").Value2 = 1.05
$ws.Range("K14").Value2 = 7.9
$ws.Range("L14").Value2 = 1.24
$ws.Range("M14").Value2 = 3.6
$ws.Range("N14").Value2 = 1.72
$ws.Range("O14").Value2 = 2
$ws.Range("P14").Value2 = 1.37
$ws.Range("Q14").Value2 = 2.85
$ws.Range("R14").Value2 = 1.75
$ws.Range("S14").Value2 = 1.98
$ws.Range("T14").Value2 = 15
$ws.Range("U14").Value2 = 30
$ws.Range("V14").Value2 = 15.5
$ws.Range("W14").Value2 = 90
$ws.Range("Z14").Value2 = 7.9
$ws.Range("AA14").Value2 = 7.3
$ws.Range("AB14").Value2 = 15
$ws.Range("AC14").Value2 = 65
$ws.Range("AD14").Value2 = 450
$ws.Range("AE14").Value2 = 7.4
$ws.Range("AF14").Value2 = 8
$ws.Range("AH14").Value2 = 12.5
$ws.Range("AI14").Value2 = 12.5
$ws.Range("AJ14").Value2 = 23
# Row 15
$ws.Range("G15").Value2 = 3.85
$ws.Range("H15").Value2 = 2.75
$ws.Range("I15").Value2 = 2.18
$ws.Range("K15").Value2 = 4.9
$ws.Range("S15").Value2 = 1.62
$ws.Range("T15").Value2 = 7.6
$ws.Range("U15").Value2 = 18.5
$ws.Range("V15").Value2 = 14
$ws.Range("W15").Value2 = 65
$ws.Range("Y15").Value2 = 70
$ws.Range("Z15").Value2 = 4.9
$ws.Range("AE15").Value2 = 5.4
$ws.Range("AF15").Value2 = 9
$ws.Range("AG15").Value2 = 9.5
$ws.Range("AH15").Value2 = 21
# Row 17
$ws.Range("G17").Value2 = 2.92
$ws.Range("H17").Value2 = 2.92
$ws.Range("I17").Value2 = 2.52
$ws.Range("J17").Value2 = 1.09
$ws.Range("K17").Value2 = 6.2
$ws.Range("L17").Value2 = 1.38
$ws.Range("M17").Value2 = 2.82
$ws.Range("T17").Value2 = 7.9
$ws.Range("U17").Value2 = 14.5
$ws.Range("V17").Value2 = 10.5
$ws.Range("W17").Value2 = 37
$ws.Range("Y17").Value2 = 37
$ws.Range("Z17").Value2 = 6.2
$ws.Range("AA17").Value2 = 5.7
$ws.Range("AE17").Value2 = 7.8
$ws.Range("AF17").Value2 = 12.5
$ws.Range("AG17").Value2 = 9.25
$ws.Range("AH17").Value2 = 29
$ws.Range("AJ17").Value2 = 30
# Row 20
$ws.Range("H20").Value2 = 3.55
$ws.Range("I20").Value2 = 2.1
$ws.Range("N20").Value2 = 1.57
$ws.Range("O20").Value2 = 2.12
$ws.Range("R20").Value2 = 1.5
$ws.Range("S20").Value2 = 2.27
$ws.Range("T20").Value2 = 13
$ws.Range("U20").Value2 = 19
$ws.Range("Y20").Value2 = 25
$ws.Range("Z20").Value2 = 14
$ws.Range("AA20").Value2 = 7.2
$ws.Range("AB20").Value2 = 11.5
$ws.Range("AC20").Value2 = 40
$ws.Range("AE20").Value2 = 10.25
$ws.Range("AF20").Value2 = 12.5
$ws.Range("AI20").Value2 = 15
$ws.Range("AJ20").Value2 = 20
# Row 22
$ws.Range("G22").Value2 = 3
$ws.Range("I22").Value2 = 2.4
$ws.Range("T22").Value2 = 9.5
$ws.Range("U22").Value2 = 15
$ws.Range("Y22").Value2 = 34
$ws.Range("Z22").Value2 = 9
$ws.Range("AE22").Value2 = 8
$ws.Range("AG22").Value2 = 9.5
# Row 23
$ws.Range("G23").Value2 = 3.6
$ws.Range("I23").Value2 = 2.25
$ws.Range("T23").Value2 = 9.5
$ws.Range("U23").Value2 = 17
$ws.Range("V23").Value2 = 13
$ws.Range("W23").Value2 = 41
$ws.Range("AE23").Value2 = 7
$ws.Range("AF23").Value2 = 10
$ws.Range("AH23").Value2 = 21
$ws.Range("AI23").Value2 = 19
# Row 24
$ws.Range("G24").Value2 = 2.4
$ws.Range("H24").Value2 = 3
$ws.Range("L24").Value2 = 1.25
$ws.Range("M24").Value2 = 3.75
$ws.Range("N24").Value2 = 1.85
$ws.Range("O24").Value2 = 1.95
$ws.Range("P24").Value2 = 1.36
$ws.Range("Q24").Value2 = 3
$ws.Range("R24").Value2 = 1.62
$ws.Range("S24").Value2 = 2.2
$ws.Range("T24").Value2 = 9.5
$ws.Range("V24").Value2 = 9.5
$ws.Range("Y24").Value2 = 26
$ws.Range("Z24").Value2 = 10
$ws.Range("AA24").Value2 = 6
$ws.Range("AB24").Value2 = 11
$ws.Range("AJ24").Value2 = 29
# Row 25
$ws.Range("G25").Value2 = 4.1
$ws.Range("H25").Value2 = 3.25
$ws.Range("I25").Value2 = 1.95
$ws.Range("AA25").Value2 = 6
$ws.Range("AG25").Value2 = 9
# Row 26
$ws.Range("J26").Value2 = 1.04
$ws.Range("K26").Value2 = 13
# Row 27
$ws.Range("G27").Value2 = 3.3
$ws.Range("H27").Value2 = 3.7
$ws.Range("I27").Value2 = 2.05
$ws.Range("AA27").Value2 = 7
$ws.Range("AF27").Value2 = 11
$ws.Range("AG27").Value2 = 9
# Row 28
$ws.Range("N28").Value2 = 1.85
$ws.Range("O28").Value2 = 1.95
# Row 29
$ws.Range("G29").Value2 = 1.48
$ws.Range("I29").Value2 = 6.5
$ws.Range("L29").Value2 = 1.2
$ws.Range("M29").Value2 = 4.33
$ws.Range("U29").Value2 = 7.5
$ws.Range("AA29").Value2 = 8.5
# Row 30
$ws.Range("G30").Value2 = 3.5
$ws.Range("H30").Value2 = 3.4
$ws.Range("I30").Value2 = 1.9
$ws.Range("K30").Value2 = 10
$ws.Range("R30").Value2 = 1.95
$ws.Range("S30").Value2 = 1.8
$ws.Range("AD30").Value2 = 351
$ws.Range("AF30").Value2 = 9
$ws.Range("AG30").Value2 = 9
# Row 35
$ws.Range("G35").Value2 = 7.5
$ws.Range("H35").Value2 = 5
$ws.Range("I35").Value2 = 1.27
$ws.Range("N35").Value2 = 1.18
$ws.Range("O35").Value2 = 4.5
$ws.Range("U35").Value2 = 51
$ws.Range("V35").Value2 = 29
$ws.Range("W35").Value2 = 101
$ws.Range("X35").Value2 = 51
$ws.Range("Y35").Value2 = 41
$ws.Range("Z35").Value2 = 34
$ws.Range("AA35").Value2 = 15
$ws.Range("AB35").Value2 = 15
$ws.Range("AC35").Value2 = 29
$ws.Range("AD35").Value2 = 67
$ws.Range("AE35").Value2 = 21
$ws.Range("AF35").Value2 = 13
$ws.Range("AG35").Value2 = 11
$ws.Range("AH35").Value2 = 13
$ws.Range("AI35").Value2 = 11
$ws.Range("AJ35").Value2 = 15
# Row 38
$ws.Range("I38").Value2 = 4.5
$ws.Range("Z38").Value2 = 17
$ws.Range("AA38").Value2 = 8.5
$ws.Range("AE38").Value2 = 17
$ws.Range("AJ38").Value2 = 34
